$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Locate "${delai} $(unite}, du ${date_debut} " and rework it into the
#    new run layout:
#       "${delai}"           (not bold)
#       " ${unite}"          (bold)
#       " "                  (bold)
#       <bookmarkStart/End name="_GoBack">
#       "du ${date_debut} "  (not bold)
# ---------------------------------------------------------------------------

$r = $d.Content
$found = $r.Find.Execute('${delai} $(unite}, du ${date_debut} ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$blockStart = $r.Start
$blockEnd = $r.End

# Clear the whole old block text first.
$whole = $d.Range($blockStart, $blockEnd)
$whole.Text = ""

$pos = $blockStart

# --- Run 1: "${delai}" -> not bold ---
$t1 = '${delai}'
$ins1 = $d.Range($pos, $pos)
$ins1.InsertBefore($t1)
$run1 = $d.Range($pos, $pos + $t1.Length)
$run1.Bold = 0
$pos = $pos + $t1.Length

# --- Run 2: " ${unite}" -> bold ---
$t2 = ' ${unite}'
$ins2 = $d.Range($pos, $pos)
$ins2.InsertBefore($t2)
$run2 = $d.Range($pos, $pos + $t2.Length)
$run2.Bold = 1
$pos = $pos + $t2.Length

# --- Run 3: " " -> bold ---
$t3 = ' '
$ins3 = $d.Range($pos, $pos)
$ins3.InsertBefore($t3)
$run3 = $d.Range($pos, $pos + $t3.Length)
$run3.Bold = 1
$pos = $pos + $t3.Length

# --- Run 4: "du ${date_debut} " -> not bold ---
$t4 = 'du ${date_debut} '
$ins4 = $d.Range($pos, $pos)
$ins4.InsertBefore($t4)
$run4 = $d.Range($pos, $pos + $t4.Length)
$run4.Bold = 0

# --- Bookmark "_GoBack" right between run 3 and run 4. Adding a bookmark
#     named "_GoBack" relocates Word's existing "_GoBack" bookmark here,
#     removing it from its previous location further down the document.
#     (Must be added *after* run 4's text is inserted so the bookmark ends
#     up right before that text rather than being pushed past it.) ---
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Host "Done"
